$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mer_indicators")

# "PMTCT FO" moves from the HS dataset (vec_mer_hs_indicators, row 28) to
# the SMI dataset (vec_mer_smi_indicators, new row next to the other SMI
# rows). First remove the old HS "PMTCT FO" row (28 is below the insertion
# point below, so removing it first keeps the row-15 index valid), then
# insert the new SMI "PMTCT FO" row at row 15, pushing rows 15-27 down to
# 16-28 and leaving the trailing rows 29-30 untouched.
$ws.Rows.Item(28).Delete()
$ws.Rows.Item(15).Insert()
$ws.Cells.Item(15, 1).Value = "vec_mer_smi_indicators"
$ws.Cells.Item(15, 3).Value = "PMTCT FO"

# Update the active selection / active sheet to match the editor's final
# on-screen state: mer_indicators becomes the active tab, with A14:C15
# selected (active cell A14).
$ws.Activate() | Out-Null
$ws.Range("A14:C15").Select() | Out-Null
